$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Structure: the cached "datetimeFigureOut" auto-date text stamped on the
#    Slide Master and every Slide Layout moved from 2017/8/28 -> 2017/9/9.
#    (ppPlaceholderDate = 16)
# ---------------------------------------------------------------------------
$ppPlaceholderDate = 16
$oldDate = "2017/8/28"
$newDate = "2017/9/9"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePh = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePh = $true
            }
        } catch {
            $isDatePh = $false
        }
        if ($isDatePh) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide Master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every Slide Layout hanging off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $layout = $layouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2) GUI adjustment: reposition the picture on slide 1 (move it down).
# ---------------------------------------------------------------------------
$emuPerPt = 12700
$s1 = $p.Slides.Item(1)
$pic1 = $s1.Shapes.Item(1)
$pic1.Top = 1492075 / $emuPerPt

# ---------------------------------------------------------------------------
# 3) New GUI widget 'choose': nudge the "图片 14" picture on slide 6.
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
for ($i = 1; $i -le $s6.Shapes.Count; $i++) {
    $shp = $s6.Shapes.Item($i)
    if ($shp.Name -eq "图片 14") {
        $shp.Left = 340877 / $emuPerPt
        $shp.Top = 3212698 / $emuPerPt
    }
}
